# SampleAttributeData.xlsx test fixture update:
# The old sample had a "Name" column (A/B/C/D) plus BRKEY/DISTRICT columns
# and an extra blank, styled column. The updated fixture drops the "Name"
# column entirely (BRKEY becomes the first column) and repurposes the
# trailing blank/styled column as a new "Inspection_Date" header, to match
# the iAM import's new inspectionDateColumnName input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete column A ("Name"/A,B,C,D) - everything else shifts one column left,
# so BRKEY (old B) becomes A, DISTRICT (old C) becomes B, and the blank
# styled column (old D) becomes C.
$ws.Columns.Item(1).Delete()

# Give that now-blank-but-styled column (C) its new header.
$ws.Range("C1").Value = "Inspection_Date"

# Match the saved selection/cursor position from the edited workbook.
$ws.Range("D4").Select()
